$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new column before column L ("P/l before tax" / "P/L Before Tax"),
# shifting existing columns L:T to M:U for all 47 rows of data.
$ws.Range("L1:L47").Insert(-4161)

# New column header text (row 1 uses the lowercase MoneyControl-style label,
# row 2 uses the Title Case label used elsewhere in the sheet).
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"
